# Weekly update: insert 4 new report rows (week of 2022-01-17, serial 44578)
# right before the existing row 359, pushing the remaining rows (old 359-382)
# down to 363-386. Dimension grows from A1:R382 to A1:R386.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 359 (shifts 359:382 -> 363:386).
$ws.Rows("359:362").Insert()

# New weekly data rows for Sandia @ Femacal de La Calera, Coquimbo.
$newRows = @(
    @(3, 'Femacal de La Calera', 'Coquimbo', 44578, 5, 100112028, 'Sandia', 'Sin especificar', 'Extra',   230, 3000, 3000, 3000, '$/unidad', 'Paine', 3000, 1, 'Hortaliza'),
    @(3, 'Femacal de La Calera', 'Coquimbo', 44578, 5, 100112028, 'Sandia', 'Sin especificar', 'Primera', 250, 2000, 2000, 2000, '$/unidad', 'Paine', 2000, 1, 'Hortaliza'),
    @(3, 'Femacal de La Calera', 'Coquimbo', 44578, 5, 100112028, 'Sandia', 'Sin especificar', 'Segunda', 280, 1500, 1500, 1500, '$/unidad', 'Paine', 1500, 1, 'Hortaliza'),
    @(3, 'Femacal de La Calera', 'Coquimbo', 44578, 5, 100112028, 'Sandia', 'Sin especificar', 'Tercera', 250, 1000, 1000, 1000, '$/unidad', 'Paine', 1000, 1, 'Hortaliza')
)

$r = 359
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
